$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Expand the "Comments ... are allowed ..." sentence in the Comment row.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "are allowed within the key, represented within regular brackets after the pipe symbol.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "are allowed within the key, represented within regular brackets after the pipe symbol. Comment can be placed both/either before and/or after key and/or value.",
    2
)

# ---------------------------------------------------------------------------
# 2) Extend the "{value|(comment) key}" example cell with two more variants:
#    {value|(comment) key} or {value (comment)|key} or {value (comment)|(comment) key}
# ---------------------------------------------------------------------------
$sel = $d.Content
$sel.Find.Execute(") key}", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insertionPoint = $sel.End
$r = $d.Range($insertionPoint, $insertionPoint)

# Insert the whole addition as plain text first ...
$addition = " or {value (comment)|key} or {value (comment)|(comment) key}"
$r.InsertAfter($addition)

# ... then re-italicise the three new "comment" occurrences inside that span.
$spanStart = $r.Start
$spanEnd = $r.End
$search = $d.Range($spanStart, $spanEnd)
$search.Find.Execute("comment", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
while ($search.Find.Found -and $search.Start -lt $spanEnd) {
    $search.Font.Italic = 1
    $search.Collapse(0)
    $search.Find.Execute("comment", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
}

# ---------------------------------------------------------------------------
# 3) Add a new bullet point to the "Open for discussion" list.
# ---------------------------------------------------------------------------
$sel2 = $d.Content
$sel2.Find.Execute(
    "Creating ontology terms from collected SOPs, and linking the keys with the ontologies.",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0
) | Out-Null
$sel2.InsertParagraphAfter()
$newPara = $sel2.Paragraphs(1).Next()
$newPara.Range.InsertBefore("How step number should be counted, e.g., should it be restarted from 1 after a new section.")
